# Updated cryptos list on Thu Apr 20 06:36:19 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns for each coin row, and
# swaps the Monero / EthereumClassic rows (27 <-> 28) back to Monero-first
# ordering. Price/volume cells hold text (not numbers) in the source data,
# so numeric-looking values are entered with a leading "'" (quote-prefix)
# to force Excel to keep them as literal text instead of auto-converting
# them to doubles (which would corrupt formatting like trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.130.79'
$ws.Range('E2').Value = '  -3.66%  '
$ws.Range('D3').Value = '1.968.74'
$ws.Range('E3').Value = '  -5.64%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = "'328.25"
$ws.Range('E5').Value = '  -3.76%  '
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('D7').Value = "'0.5007"
$ws.Range('E7').Value = '  -5.19%  '
$ws.Range('D8').Value = "'0.4222"
$ws.Range('E8').Value = '  -3.67%  '
$ws.Range('D9').Value = "'52.74"
$ws.Range('E9').Value = '  -3.75%  '
$ws.Range('D10').Value = "'0.09196"
$ws.Range('E10').Value = '  -1.57%  '
$ws.Range('D11').Value = "'1.103"
$ws.Range('E11').Value = '  -6.01%  '
$ws.Range('D12').Value = "'23.09"
$ws.Range('E12').Value = '  -5.77%  '
$ws.Range('D13').Value = '2.013.49'
$ws.Range('E13').Value = '  -2.87%  '
$ws.Range('D14').Value = "'7.906"
$ws.Range('E14').Value = '  -6.93%  '
$ws.Range('D15').Value = "'6.450"
$ws.Range('E15').Value = '  -5.99%  '
$ws.Range('D16').Value = "'1.008"
$ws.Range('E16').Value = '  +0.29%  '
$ws.Range('D17').Value = "'0.00001105"
$ws.Range('E17').Value = '  -4.47%  '
$ws.Range('D18').Value = "'91.69"
$ws.Range('E18').Value = '  -9.66%  '
$ws.Range('D19').Value = "'0.06709"
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('D20').Value = "'19.36"
$ws.Range('E20').Value = '  -7.95%  '
$ws.Range('D21').Value = "'1.004"
$ws.Range('E21').Value = '  +0.16%  '
$ws.Range('E22').Value = '  -4.62%  '
$ws.Range('D23').Value = '29.173.18'
$ws.Range('E23').Value = '  -3.60%  '
$ws.Range('D24').Value = "'12.13"
$ws.Range('E24').Value = '  -2.07%  '
$ws.Range('D25').Value = "'2.286"
$ws.Range('E25').Value = '  -1.56%  '
$ws.Range('D26').Value = '2.265.33'
$ws.Range('E26').Value = '  -2.95%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Value = "'156.72"
$ws.Range('E27').Value = '  -3.46%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = "'20.64"
$ws.Range('E28').Value = '  -5.20%  '
$ws.Range('D29').Value = "'6.206"
$ws.Range('E29').Value = '  -9.20%  '
$ws.Range('D30').Value = "'2.269"
$ws.Range('E30').Value = '  -8.50%  '
$ws.Range('D31').Value = "'126.67"
$ws.Range('E31').Value = '  -5.08%  '
$ws.Range('D32').Value = "'1.047"
$ws.Range('E32').Value = '  -7.05%  '
$ws.Range('D33').Value = "'0.09866"
$ws.Range('E33').Value = '  -5.79%  '
$ws.Range('D34').Value = "'1.536"
$ws.Range('E34').Value = '  -7.46%  '
$ws.Range('D35').Value = "'5.793"
$ws.Range('E35').Value = '  -7.30%  '
$ws.Range('D36').Value = "'3.682"
$ws.Range('E36').Value = '  -5.92%  '
$ws.Range('D37').Value = "'0.02433"
$ws.Range('E37').Value = '  -6.69%  '
$ws.Range('D38').Value = "'9.053"
$ws.Range('E38').Value = '  -8.61%  '
$ws.Range('D39').Value = "'1.301"
$ws.Range('E39').Value = '  -3.01%  '
$ws.Range('D40').Value = "'0.06370"
$ws.Range('E40').Value = '  -5.41%  '
$ws.Range('D41').Value = "'0.6468"
$ws.Range('E41').Value = '  -6.66%  '
$ws.Range('D42').Value = "'11.46"
$ws.Range('E42').Value = '  -8.80%  '
$ws.Range('D43').Value = "'0.1989"
$ws.Range('E43').Value = '  -9.48%  '
$ws.Range('E44').Value = '  +0.17%  '
$ws.Range('D45').Value = "'0.6257"
$ws.Range('E45').Value = '  -7.07%  '
$ws.Range('D46').Value = "'13.44"
$ws.Range('E46').Value = '  -5.35%  '
$ws.Range('D47').Value = "'2.202"
$ws.Range('E47').Value = '  -7.45%  '
$ws.Range('D48').Value = "'1.284"
$ws.Range('E48').Value = '  -0.43%  '
$ws.Range('D49').Value = "'3.474"
$ws.Range('E49').Value = '  -4.42%  '
$ws.Range('D50').Value = "'0.00000000331"
$ws.Range('E50').Value = '  -4.98%  '
$ws.Range('D51').Value = "'0.06987"
